$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.680.08'
$ws.Range("E2").Value = '  +2.91%  '
$ws.Range("D3").Value = '2.199.16'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'258.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.72%  '
$ws.Range("D6").Value = "'83.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +11.18%  '
$ws.Range("D7").Value = "'0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.00%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = "'0.597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = "'44.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.85%  '
$ws.Range("D11").Value = "'0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = "'7.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.62%  '
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = '2.525.94'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = "'14.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '2.197.20'
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = "'0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").Value = '43.575.37'
$ws.Range("E18").Value = '  +2.82%  '
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("D20").Value = "'69.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").Value = "'5.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.76%  '
$ws.Range("E22").Value = '  +11.72%  '
$ws.Range("D23").Value = "'231.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.05%  '
$ws.Range("D24").Value = "'8.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.90%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").Value = "'10.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").Value = "'39.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.89%  '
$ws.Range("E29").Value = '  +2.97%  '
$ws.Range("D30").Value = "'2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.68%  '
$ws.Range("D31").Value = "'174.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("D32").Value = "'20.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").Value = "'0.0859"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("E34").Value = '  +3.79%  '
$ws.Range("E35").Value = '  +2.02%  '
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("D37").Value = "'4.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.24%  '
$ws.Range("D38").Value = "'0.0360"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("D39").Value = "'12.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.58%  '
$ws.Range("E40").Value = '  +8.30%  '
$ws.Range("D41").Value = "'2.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.57%  '
$ws.Range("D42").Value = "'62.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.01%  '
$ws.Range("D43").Value = "'5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.12%  '
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("D45").Value = "'0.0978"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("D46").Value = "'8.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("D47").Value = "'99.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("E48").Value = '  +5.90%  '
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").Value = "'0.437"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.04%  '
$ws.Range("E51").Value = '  +7.39%  '
